# Reatraining the quarterly forecast model for PCSunEnergy
# Shift the Notified Production Wind series forward by 3 days (new data
# pull) and refresh the forecast values in column B for rows 2..97.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(531,535,542,548,617,621,615,617,736,730,738,750,1072,1083,1100,1116,1602,1615,1629,1639,1948,1965,1963,1969,2042,2070,2068,2057,2029,2031,2035,2033,1905,1902,1899,1896,1774,1770,1767,1763,1595,1591,1586,1580,1333,1326,1320,1313,1018,1013,1007,1000,815,810,806,801,706,703,700,695,617,615,613,610,597,595,593,591,561,559,557,555,553,551,549,548,576,576,576,577,675,675,675,676,752,752,753,754,781,781,782,783,0,0,0,0)

for ($row = 2; $row -le 97; $row++) {
    $idx = $row - 2
    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.Value = $dateCell.Value2 + 3
    $ws.Cells.Item($row, 2).Value = $newValues[$idx]
}
